$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new "map display" issue row (row 9), copying the style
# already used by the other issue rows (e.g. row 8: Sprint/title/desc/status).
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A9").Value = "Sprint 3"
$ws.Range("C9").Value = "In 2D mode, the map can be dragged untill seeing the eadges of the map "
$ws.Range("B9").Value = "map display"
$ws.Range("D9").Value = "Open"

$ws.Rows.Item(9).RowHeight = 45

# Update the view to match: scrolled down to row 7, with C9 selected.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C9").Select()
